$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Columns.Item(2).ColumnWidth = 15.035714285714286
$ws.Columns.Item(4).ColumnWidth = 14.535714285714286
$ws.Columns.Item(5).ColumnWidth = 15.910714285714286
$ws.Columns.Item(6).ColumnWidth = 15.410714285714286
$ws.Columns.Item(7).ColumnWidth = 17.285714285714285
$ws.Columns.Item(8).ColumnWidth = 16.785714285714285
